$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2863.9119
$ws.Range("I62").Value = 2379.24
$ws.Range("J62").Value = 4210.222
$ws.Range("K62").Value = 2379.24
$ws.Range("L62").Value = 4210.222
$ws.Range("M62").Value = -1755.24
$ws.Range("N62").Value = -5458.222

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 2863.9119
$ws.Range("I65").Value = 2379.24
$ws.Range("J65").Value = 4210.222
$ws.Range("K65").Value = 11896.2
$ws.Range("L65").Value = 21051.11
$ws.Range("M65").Value = -8776.199999999999
$ws.Range("N65").Value = -27291.11

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 699.38464
$ws.Range("I98").Value = 399.2
$ws.Range("K98").Value = 399.2
$ws.Range("M98").Value = 1098.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 500080
$ws.Range("I103").Value = 500080
$ws.Range("K103").Value = 1500240
$ws.Range("M103").Value = -1499654

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 9525675
$ws.Range("I106").Value = 13334505
$ws.Range("K106").Value = 13334505
$ws.Range("M106").Value = -13333874

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 397.52173
$ws.Range("I107").Value = 497.52942
$ws.Range("J107").Value = 114.166664
$ws.Range("K107").Value = 497.52942
$ws.Range("L107").Value = 114.166664
$ws.Range("M107").Value = 1422.47058
$ws.Range("N107").Value = -3954.166664

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 699.38464
$ws.Range("I122").Value = 399.2
$ws.Range("K122").Value = 1197.6
$ws.Range("M122").Value = 1252.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 387331.38
$ws.Range("J129").Value = 419567.34
$ws.Range("L129").Value = 1258702.02
$ws.Range("N129").Value = -1268702.02

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2850.5667
$ws.Range("I132").Value = 3041.6785
$ws.Range("J132").Value = 175
$ws.Range("K132").Value = 9125.0355
$ws.Range("L132").Value = 525
$ws.Range("M132").Value = -6595.0355
$ws.Range("N132").Value = -5585

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 102580.1
$ws.Range("I137").Value = 2500
$ws.Range("K137").Value = 7500
$ws.Range("M137").Value = -4950

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1688.7
$ws.Range("I138").Value = 454.6
$ws.Range("J138").Value = 2511.4333
$ws.Range("K138").Value = 1363.8
$ws.Range("L138").Value = 7534.2999
$ws.Range("M138").Value = 3776.2
$ws.Range("N138").Value = -17814.2999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 17550.572
$ws.Range("J24").Value = 17550.572
$ws.Range("L24").Value = 17550.572
$ws.Range("N24").Value = -18298.572

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17345.975
$ws.Range("I32").Value = 18597.29
$ws.Range("K32").Value = 18597.29
$ws.Range("M32").Value = -18310.29

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 24391918
$ws.Range("I74").Value = 27779596
$ws.Range("J74").Value = 635.8
$ws.Range("K74").Value = 27779596
$ws.Range("L74").Value = 635.8
$ws.Range("M74").Value = -27778722
$ws.Range("N74").Value = -2383.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 24391918
$ws.Range("I77").Value = 27779596
$ws.Range("J77").Value = 635.8
$ws.Range("K77").Value = 138897980
$ws.Range("L77").Value = 3179
$ws.Range("M77").Value = -138893612
$ws.Range("N77").Value = -11915

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H100").Value = 17550.572
$ws.Range("J100").Value = 17550.572
$ws.Range("L100").Value = 17550.572
$ws.Range("N100").Value = -19714.572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4200
$ws.Range("I105").Value = 4533.3335
$ws.Range("J105").Value = 3533.3333
$ws.Range("K105").Value = 4533.3335
$ws.Range("L105").Value = 3533.3333
$ws.Range("M105").Value = -2786.3335
$ws.Range("N105").Value = -7027.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 31892.514
$ws.Range("I134").Value = 38487.434
$ws.Range("K134").Value = 115462.302
$ws.Range("M134").Value = -112927.302

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9599.465
$ws.Range("I31").Value = 11350.871
$ws.Range("J31").Value = 5075
$ws.Range("K31").Value = 11350.871
$ws.Range("L31").Value = 5075
$ws.Range("M31").Value = -11055.871
$ws.Range("N31").Value = -5665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 9599.465
$ws.Range("I34").Value = 11350.871
$ws.Range("J34").Value = 5075
$ws.Range("K34").Value = 11350.871
$ws.Range("L34").Value = 5075
$ws.Range("M34").Value = -11148.871
$ws.Range("N34").Value = -5479

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 815.41
$ws.Range("I131").Value = 282.5
$ws.Range("J131").Value = 837.61456
$ws.Range("K131").Value = 847.5
$ws.Range("L131").Value = 2512.84368
$ws.Range("M131").Value = 4192.5
$ws.Range("N131").Value = -12592.84368

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 2195.6316
$ws.Range("I136").Value = 909.0769
$ws.Range("J136").Value = 4983.1665
$ws.Range("K136").Value = 2727.2307
$ws.Range("L136").Value = 14949.4995
$ws.Range("M136").Value = 2372.7693
$ws.Range("N136").Value = -25149.4995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 1991.6666
$ws.Range("I141").Value = 1675.1052
$ws.Range("K141").Value = 5025.3156
$ws.Range("M141").Value = 154.6844000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 18369.23
$ws.Range("I70").Value = 19133.334
$ws.Range("K70").Value = 19133.334
$ws.Range("M70").Value = -18863.334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 18369.23
$ws.Range("I73").Value = 19133.334
$ws.Range("K73").Value = 19133.334
$ws.Range("M73").Value = -18197.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1895.7142
$ws.Range("I22").Value = 1404.8
$ws.Range("K22").Value = 1404.8
$ws.Range("M22").Value = -1109.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1895.7142
$ws.Range("I27").Value = 1404.8
$ws.Range("K27").Value = 1404.8
$ws.Range("M27").Value = -1297.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1954.4546
$ws.Range("I68").Value = 1886.7333
$ws.Range("J68").Value = 2099.5715
$ws.Range("K68").Value = 1886.7333
$ws.Range("L68").Value = 2099.5715
$ws.Range("M68").Value = -1137.7333
$ws.Range("N68").Value = -3597.5715

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1954.4546
$ws.Range("I71").Value = 1886.7333
$ws.Range("J71").Value = 2099.5715
$ws.Range("K71").Value = 9433.666500000001
$ws.Range("L71").Value = 10497.8575
$ws.Range("M71").Value = -5689.666500000001
$ws.Range("N71").Value = -17985.8575

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 1310476.2
$ws.Range("I122").Value = 2803963.5
$ws.Range("J122").Value = 3674.875
$ws.Range("K122").Value = 8411890.5
$ws.Range("L122").Value = 11024.625
$ws.Range("M122").Value = -8409440.5
$ws.Range("N122").Value = -15924.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1366.6666
$ws.Range("I132").Value = 1250
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 3750
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -1220
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4560.4
$ws.Range("I62").Value = 3901
$ws.Range("K62").Value = 3901
$ws.Range("M62").Value = -3277

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 4560.4
$ws.Range("I65").Value = 3901
$ws.Range("K65").Value = 19505
$ws.Range("M65").Value = -16385
